$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: collapse the 9-run title "Преобразование «ты - высказываний»
# в развернутые «я - высказывания»." into one contiguous run with the
# same full text (removes the run fragmentation / extra spacing runs).
# ---------------------------------------------------------------------
$oldTitle = "Преобразование «ты - высказываний» в развернутые «я - высказывания»."
$newTitle = "Преобразование «ты - высказываний» в развернутые «я - высказывания»."
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: replace the last paragraph (which currently just holds a
# lone <w:tab/> plus the trailing _GoBack bookmark) with:
#   - the same tab run, now followed by a new explanatory-text run;
#   - an empty spacer paragraph;
#   - a new numbered item "«Не серди меня»" + its explanatory paragraph;
#   - another empty spacer paragraph;
#   - a new numbered item "«Ну ты даешь»" + its (multi-run) explanatory
#     paragraph;
#   - the _GoBack bookmark, now re-anchored at the very end.
# Replacing the *entire* paragraph range (through its end-of-paragraph
# mark) in one InsertXML call lets us drop the old bookmark location
# and re-place it at the true end, instead of having it get stranded
# mid-document.
# ---------------------------------------------------------------------
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range

$replacementXml = @"
<w:p $wordNs w14:paraId="4F3D9B62" w14:textId="097ECF15" w:rsidR="001E4270" w:rsidRPr="001E4270" w:rsidRDefault="001E4270" w:rsidP="001E4270">
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t>Мне очень жаль, что ты поверил этому слуху. Я хотела бы чтобы в наших отношениях было больше доверия. Давай договоримся, что будем обсуждать все противоречивые ситуации.</w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">«Не серди меня» </w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:tab/>
    <w:t>Мы с тобой уже неоднократно обсуждали эту ситуацию, но сегодня, когда я пришел домой, я снова увидел гору грязной посуды. Меня огорчает, что наши договоренности не соблюдаются. Мне хотелось бы быть услышанным.</w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>«Ну ты даешь»</w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="993"/>
    </w:tabs>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">Я попросил тебя помочь мне выполнить задание, </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">ведь ты гораздо лучше разбираешься в этой теме, </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">но ты забыла о </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="ru-RU"/>
    </w:rPr>
    <w:t>нем. Мне очень нужна твоя помощь и мне обидно, что ты недооцениваешь важность моей просьбы. Давай попробуем помогать друг другу больше.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@

$lastRange.InsertXML($replacementXml)
